$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 55 entirely: rows 56:98 shift up to become 55:97
$ws.Rows.Item(55).Delete()

# Scroll the window so row 46 is at the top, matching the post-edit view
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 4

# Select the row that now occupies the deleted row's place (full row 55),
# with the active cell over column D
$ws.Range("D55").Select()
$ws.Rows.Item(55).Select()
